$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C16").Value = "1143326713"
$ws.Range("D16").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 46667
$ws.Range("G16").Value = 1400000

$ws.Range("C17").Value = "1143326713"
$ws.Range("D17").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E17").Value = "1902"
$ws.Range("F17").Value = 56000
$ws.Range("G17").Value = 1400000

$ws.Range("C18").Value = "1143326713"
$ws.Range("D18").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E18").Value = "1901"
$ws.Range("F18").Value = 56000
$ws.Range("G18").Value = 1400000

$ws.Range("C19").Value = "1143326713"
$ws.Range("D19").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E19").Value = "1812"
$ws.Range("F19").Value = 56000
$ws.Range("G19").Value = 1400000

$ws.Range("C20").Value = "1143326713"
$ws.Range("D20").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E20").Value = "1811"
$ws.Range("F20").Value = 56000
$ws.Range("G20").Value = 1400000

$ws.Range("C21").Value = "1143326713"
$ws.Range("D21").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E21").Value = "1810"
$ws.Range("F21").Value = 56000
$ws.Range("G21").Value = 1400000

$ws.Range("C22").Value = "1143326713"
$ws.Range("D22").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E22").Value = "1809"
$ws.Range("F22").Value = 56000
$ws.Range("G22").Value = 1400000

$ws.Range("C23").Value = "1143326713"
$ws.Range("D23").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E23").Value = "1808"
$ws.Range("F23").Value = 56000
$ws.Range("G23").Value = 1400000

$ws.Range("C24").Value = "1143326713"
$ws.Range("D24").Value = "NATALI CONTRERAS SOLAR"
$ws.Range("E24").Value = "1807"
$ws.Range("F24").Value = 37334
$ws.Range("G24").Value = 1400000

$ws.Range("C25").Value = "1050969563"
$ws.Range("D25").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E25").Value = "1903"
$ws.Range("F25").Value = 26041
$ws.Range("G25").Value = 781242

$ws.Range("C26").Value = "1050969563"
$ws.Range("D26").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E26").Value = "1902"
$ws.Range("F26").Value = 31249
$ws.Range("G26").Value = 781242

$ws.Range("C27").Value = "1050969563"
$ws.Range("D27").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E27").Value = "1901"
$ws.Range("F27").Value = 31249
$ws.Range("G27").Value = 781242

$ws.Range("C28").Value = "1050969563"
$ws.Range("D28").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E28").Value = "1812"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242

$ws.Range("C29").Value = "1050969563"
$ws.Range("D29").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E29").Value = "1811"
$ws.Range("F29").Value = 31249
$ws.Range("G29").Value = 781242

$ws.Range("C30").Value = "1050969563"
$ws.Range("D30").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E30").Value = "1810"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242

$ws.Range("C31").Value = "1050969563"
$ws.Range("D31").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E31").Value = "1809"
$ws.Range("F31").Value = 31249
$ws.Range("G31").Value = 781242

$ws.Range("C32").Value = "1050969563"
$ws.Range("D32").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E32").Value = "1808"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242

$ws.Range("C33").Value = "1050969563"
$ws.Range("D33").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E33").Value = "1807"
$ws.Range("F33").Value = 31249
$ws.Range("G33").Value = 781242

$ws.Range("C34").Value = "1050969563"
$ws.Range("D34").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E34").Value = "1806"
$ws.Range("F34").Value = 31249
$ws.Range("G34").Value = 781242

$ws.Range("C35").Value = "1050969563"
$ws.Range("D35").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E35").Value = "1805"
$ws.Range("F35").Value = 31249
$ws.Range("G35").Value = 781242

$ws.Range("C36").Value = "1050969563"
$ws.Range("D36").Value = "CAMILO ANDRES DIAZ PINO"
$ws.Range("E36").Value = "1804"
$ws.Range("F36").Value = 22916
$ws.Range("G36").Value = 781242

$ws.Range("C37").Value = "3815072"
$ws.Range("D37").Value = "HUGO RAFAEL CERA TORRES"
$ws.Range("E37").Value = "1807"
$ws.Range("F37").Value = 20833
$ws.Range("G37").Value = 781242

$ws.Range("C38").Value = "3976400"
$ws.Range("D38").Value = "EFRAIN POSSO AYALA"
$ws.Range("E38").Value = "1903"
$ws.Range("F38").Value = 26041
$ws.Range("G38").Value = 781242

$ws.Range("C39").Value = "3976400"
$ws.Range("D39").Value = "EFRAIN POSSO AYALA"
$ws.Range("E39").Value = "1902"
$ws.Range("F39").Value = 31249
$ws.Range("G39").Value = 781242

$ws.Range("C40").Value = "3976400"
$ws.Range("D40").Value = "EFRAIN POSSO AYALA"
$ws.Range("E40").Value = "1901"
$ws.Range("F40").Value = 31249
$ws.Range("G40").Value = 781242

$ws.Range("C41").Value = "3976400"
$ws.Range("D41").Value = "EFRAIN POSSO AYALA"
$ws.Range("E41").Value = "1812"
$ws.Range("F41").Value = 31249
$ws.Range("G41").Value = 781242

$ws.Range("C42").Value = "3976400"
$ws.Range("D42").Value = "EFRAIN POSSO AYALA"
$ws.Range("E42").Value = "1811"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242

$ws.Range("C43").Value = "3976400"
$ws.Range("D43").Value = "EFRAIN POSSO AYALA"
$ws.Range("E43").Value = "1810"
$ws.Range("F43").Value = 31249
$ws.Range("G43").Value = 781242

$ws.Range("C44").Value = "3976400"
$ws.Range("D44").Value = "EFRAIN POSSO AYALA"
$ws.Range("E44").Value = "1809"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242

$ws.Range("C45").Value = "3976400"
$ws.Range("D45").Value = "EFRAIN POSSO AYALA"
$ws.Range("E45").Value = "1808"
$ws.Range("F45").Value = 31249
$ws.Range("G45").Value = 781242

$ws.Range("C46").Value = "3976400"
$ws.Range("D46").Value = "EFRAIN POSSO AYALA"
$ws.Range("E46").Value = "1807"
$ws.Range("F46").Value = 15625
$ws.Range("G46").Value = 781242
